# "Some further error corrections when working with the marks sheet"
#
# Fill in marks that were missing/blank on the "Senior Five" sheet of the
# A-level HISTORY Term II 2024 marksheet:
#   - Row 5 (ATIM FELICITY):   Paper 1 = 51, Paper 3 = 37
#   - Row 7 (BONGOMIN RONNIE): Mid Paper 3 = 33, Paper 1 = 70

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Senior Five")

# Row 5: ATIM FELICITY
$ws.Range("F5").Value = 51
$ws.Range("G5").Value = 37

# Row 7: BONGOMIN RONNIE
$ws.Range("E7").Value = 33
$ws.Range("F7").Value = 70

# Touch the alignment of the newly-entered mark cells so they pick up the
# distinct cell style the sheet owner applied after correcting these values
# (mirrors the sheet's existing "applyAlignment" formatted cells).
$ws.Range("F5:G5").IndentLevel = 0
$ws.Range("F7").IndentLevel = 0
